$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values to set for rows 2 and 3 (same values in both rows).
$values = @{
    "D" = 0.209
    "E" = 0.32
    "F" = 0.237
    "G" = 0.1729173867957138
    "H" = 0.1271171793985482
    "I" = 0.567749740753543
    "J" = 0.4655000457073643
    "K" = 476.2
    "L" = 0.4115105426892499
    "M" = 238.698
    "N" = 0.1612824324324324
    "O" = 0.5012557748845023
    "P" = 238.698
    "Q" = 0.1612824324324324
    "R" = 0.5012557748845023
    "W" = 1.332773579624965
    "X" = 0.02919174153823396
    "Y" = 1.303581838086731
    "Z" = 3.065430463576159
    "AA" = 1.426958020907449
    "AB" = 0.02909967244874467
    "AC" = 1.397858348458704
    "AD" = 20.6
    "AE" = 0
    "AF" = 20.6
    "AG" = 20.6
    "AH" = 0.01372784219645475
    "AI" = 0.02604298356510746
    "AJ" = 0.01372784219645475
    "AK" = 0.02604298356510746
    "AL" = 67.90000000000001
    "AM" = 67.90000000000001
    "AN" = 0.03084294055996407
    "AO" = 9.675994108983799
    "AP" = 0.03084294055996407
    "AQ" = 9.675994108983799
}

foreach ($row in 2, 3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
